# "better text sizing on fig 2 go terms in silico"
# Populate the (previously empty) "for python" sheet with the condensed
# GO-CC summary table, give column A a wider custom width, and switch the
# active-tab/selection state from "sorted and combined" to "for python".

$wb = $excel.ActiveWorkbook

$wsData   = $wb.Worksheets.Item("sorted and combined")
$wsTarget = $wb.Worksheets.Item("for python")

# Header row
$headers = @("GO term", "Day 0 TD", "Day 2 TD", "Day 5 TD", "Day 12 TD", "Day 0 ND", "Day 2 ND", "Day 5 ND", "Day 12 ND")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $wsTarget.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Condensed category rows (label, Day0 TD, Day2 TD, Day5 TD, Day12 TD, Day0 ND, Day2 ND, Day5 ND, Day12 ND)
$rows = @(
    @("Chloroplast",             227, 265, 143, 135,  98, 217,  88, 60),
    @("Photosystem I",            15,  15,  15,  15,   0,   7,   6,  8),
    @("Photosystem II",           28,  34,  24,  22,   2,  32,  11,  3),
    @("Membrane",                256, 202, 127, 105,  82, 133,  43, 25),
    @("Mitochondria",             34,  76,   7,   0,  12,  35,   0,  0),
    @("Endoplasmic reticulum",    42,  29,   1,   0,   0,  16,   0,  0),
    @("Golgi appartatus",         29,   1,   0,   0,   0,  14,   0,  0),
    @("Ribosome",                 64,  94,   9,   0,  27,  41,   0,  0),
    @("Nucleus",                  19,  35,   6,   0,  10,  16,   0,  0),
    @("Cytoplasm",                136, 213,  8,   0,  39,  76,  11,  0),
    @("Secretory",                 87, 124,  4,   9,   3,  47,   0,  0)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $wsTarget.Cells.Item($r + 2, 1).Value = $row[0]
    for ($c = 1; $c -lt $row.Length; $c++) {
        $wsTarget.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Widen column A to fit the longer labels
$wsTarget.Columns.Item(1).ColumnWidth = 20

# Move the active selection on the old sheet off of A73 and onto I10
[void]$wsData.Range("I10").Select()

# Make "for python" the active/selected sheet, with A15 selected
[void]$wsTarget.Range("A15").Select()
[void]$wsTarget.Activate()
